# Generate Report for Handback
# Updates the existing handback row (new GUID/timestamps for the file that
# was already reported) and appends a brand-new handback row for a second
# file, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$guid1 = "8315e09e-2af9-4327-a12d-ac5760e73a7f"
$guid2 = "b9747646-22d8-4a72-b1f8-4868631950d6"
$hash1 = "7315c86a1ee601c1b66d28aeae6ec5dabb543a56"
$hash2 = "b1a93591cbe06d96c1ac56d822548f6042eaa4e5"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

function Set-Text($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-DateText($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = $dateFmt
}

function Add-Link($ws, $addr, $url, $display) {
    $ws.Range($addr).Value = $display
    $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $display) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Hyperlinks.Delete() on this runtime clears every hyperlink on the sheet,
# so wipe them up-front and re-add all of them (existing + new) in order.
$wsOv.Hyperlinks.Delete()

# -- Row 2 (existing file) gets the refreshed GUID / timestamp
Set-Text     $wsOv "A2" "$guid1.md"
Add-Link     $wsOv "B2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeca6f86891acf83507447d735b19c551748a951/e2e/$guid1.md" "e2e\$guid1.md"
Set-Text     $wsOv "C2" ".md"
Set-Text     $wsOv "E2" "Handed back: in sync with en-US"
Set-Text     $wsOv "F2" "Handed back: in sync with en-US"
Set-DateText $wsOv "G2" "2016-08-16 18:55:52"

# -- Row 3 (new file) is appended
Set-Text     $wsOv "A3" "$guid2.md"
Add-Link     $wsOv "B3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeca6f86891acf83507447d735b19c551748a951/e2e/$guid2.md" "e2e\$guid2.md"
Set-Text     $wsOv "C3" ".md"
Set-Text     $wsOv "E3" "Handed back: in sync with en-US"
Set-Text     $wsOv "F3" "Handed back: in sync with en-US"
Set-DateText $wsOv "G3" "2016-08-16 18:55:52"

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

# -- Row 2 (existing file)
Add-Link     $wsZh "A2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeca6f86891acf83507447d735b19c551748a951/e2e/$guid1.md" "$guid1.md"
Set-Text     $wsZh "B2" ".md"
Set-Text     $wsZh "C2" "Handed back: in sync with en-US"
Set-Text     $wsZh "D2" "e2e"
Set-Text     $wsZh "E2" "ht"
Set-Text     $wsZh "F2" "False"
Set-Text     $wsZh "G2" "$guid1.$hash1.zh-cn.xlf"
Set-DateText $wsZh "H2" "2016-08-16 18:55:47"
Add-Link     $wsZh "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ed9c988e68fc48378d6223999530836f0ffb8114/e2e/$guid1.md" "$guid1.md"
Set-Text     $wsZh "J2" "$guid1.$hash1.zh-cn.xlf"
Set-DateText $wsZh "K2" "2016-08-16 18:56:21"
Set-Text     $wsZh "M2" "True"
Set-Text     $wsZh "O2" "False"

# -- Row 3 (new file)
Add-Link     $wsZh "A3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeca6f86891acf83507447d735b19c551748a951/e2e/$guid2.md" "$guid2.md"
Set-Text     $wsZh "B3" ".md"
Set-Text     $wsZh "C3" "Handed back: in sync with en-US"
Set-Text     $wsZh "D3" "e2e"
Set-Text     $wsZh "E3" "ht"
Set-Text     $wsZh "F3" "True"
Set-Text     $wsZh "G3" "$guid2.$hash2.zh-cn.xlf"
Set-DateText $wsZh "H3" "2016-08-16 18:55:47"
Add-Link     $wsZh "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ed9c988e68fc48378d6223999530836f0ffb8114/e2e/$guid2.md" "$guid2.md"
Set-Text     $wsZh "J3" "$guid2.$hash2.zh-cn.xlf"
Set-DateText $wsZh "K3" "2016-08-16 18:56:21"
Set-Text     $wsZh "M3" "True"
Set-Text     $wsZh "O3" "False"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

# -- Row 2 (existing file)
Add-Link     $wsDe "A2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeca6f86891acf83507447d735b19c551748a951/e2e/$guid1.md" "$guid1.md"
Set-Text     $wsDe "B2" ".md"
Set-Text     $wsDe "C2" "Handed back: in sync with en-US"
Set-Text     $wsDe "D2" "e2e"
Set-Text     $wsDe "E2" "ht"
Set-Text     $wsDe "F2" "False"
Set-Text     $wsDe "G2" "$guid1.$hash1.de-de.xlf"
Set-DateText $wsDe "H2" "2016-08-16 18:55:52"
Add-Link     $wsDe "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1ed44f89aef68c44ce52bf11a7668e49f86821ee/e2e/$guid1.md" "$guid1.md"
Set-Text     $wsDe "J2" "$guid1.$hash1.de-de.xlf"
Set-DateText $wsDe "K2" "2016-08-16 18:56:28"
Set-Text     $wsDe "M2" "True"
Set-Text     $wsDe "O2" "False"

# -- Row 3 (new file)
Add-Link     $wsDe "A3" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeca6f86891acf83507447d735b19c551748a951/e2e/$guid2.md" "$guid2.md"
Set-Text     $wsDe "B3" ".md"
Set-Text     $wsDe "C3" "Handed back: in sync with en-US"
Set-Text     $wsDe "D3" "e2e"
Set-Text     $wsDe "E3" "ht"
Set-Text     $wsDe "F3" "True"
Set-Text     $wsDe "G3" "$guid2.$hash2.de-de.xlf"
Set-DateText $wsDe "H3" "2016-08-16 18:55:52"
Add-Link     $wsDe "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1ed44f89aef68c44ce52bf11a7668e49f86821ee/e2e/$guid2.md" "$guid2.md"
Set-Text     $wsDe "J3" "$guid2.$hash2.de-de.xlf"
Set-DateText $wsDe "K3" "2016-08-16 18:56:28"
Set-Text     $wsDe "M3" "True"
Set-Text     $wsDe "O3" "False"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
